# Update the weekly schedule grid on the "schedule_data" sheet.
# Row 1 = headers (days of week); Rows 2-5 = periods 1-4; Columns B..G = Mon..Sat.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule_data")

# Write the new subject names in the same order they were typed in so the
# shared-string table is rebuilt with a matching sequence.
$ws.Range("B2").Value = "情報理論"
$ws.Range("C3").Value = "統計学"
$ws.Range("C4").Value = "社会"
$ws.Range("D4").Value = "英語"
$ws.Range("D2").Value = "算数"
$ws.Range("E3").Value = "実習"
$ws.Range("F4").Value = "ドイツ語"
$ws.Range("F5").Value = "物理"
$ws.Range("E4").Value = "実習"

# Clear cells whose subjects were removed from the schedule.
$ws.Range("E2").Value = ""
$ws.Range("B3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("C5").Value = ""

# Update the active selection as recorded in the workbook (F7).
$ws.Range("F7").Select()
